$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "CEMT-class" column (D) values keyed by row number
$values = @{
    1  = "CEMT-class"
    2  = 0
    3  = "I"
    4  = "II"
    5  = "III"
    6  = "III"
    7  = "III"
    8  = "IVa"
    9  = "IVa"
    10 = "Va"
    11 = "Va"
    12 = "VIa"
    13 = "VIa"
    14 = "VIa"
    15 = "I"
    16 = "I"
    17 = "IVb"
    18 = "Vb"
    19 = "VIa"
    20 = "VIa"
    21 = "VIb"
    22 = "I"
    23 = "II"
    24 = "-"
    25 = "III"
    26 = "IV"
    27 = "Va"
    28 = "Va"
    29 = "Va"
    30 = "Vb"
    31 = "VIa"
    32 = "VIb"
    33 = "VIIa"
    34 = "VIc"
}

for ($r = 1; $r -le 34; $r++) {
    $ws.Cells.Item($r, 4).Value = $values[$r]
}

# Header formatting: bold, like the other header cells
$ws.Range("D1").Font.Bold = $true

# D2 value is numeric and left-aligned
$ws.Range("D2").HorizontalAlignment = -4131

# Update selection to reflect the new data range that was worked on
$ws.Range("A24:D34").Select() | Out-Null
